$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Append a new row of data to the bottom of the log:
#   Date 11/14/2025 (serial 45975), Error Count 14
# Copy the date formatting from the row above (A12) so the new date cell
# picks up the same existing "short date" style instead of creating a new one.
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A13").Value = 45975
$ws.Range("B13").Value = 14

# Leave the selection where the author apparently clicked next when done.
$ws.Range("C17").Select()
